$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-04-28 Sunday" "2024-04-29 Monday"

Replace-Text "979×6=" "822×2="
Replace-Text "884×4=" "669×7="
Replace-Text "773×3=" "850×2="
Replace-Text "814×2=" "983×8="
Replace-Text "818×5=" "143×4="

Replace-Text "939×3=" "139×4="
Replace-Text "769×2=" "213×5="
Replace-Text "419×9=" "627×8="
Replace-Text "471×8=" "307×2="
Replace-Text "640×5=" "586×7="

Replace-Text "757×7=" "518×6="
Replace-Text "709×6=" "371×6="
Replace-Text "293×4=" "502×3="
Replace-Text "175×3=" "268×8="
Replace-Text "236×2=" "506×3="

Replace-Text "850×7=" "245×4="
Replace-Text "766×3=" "954×2="
Replace-Text "297×5=" "988×9="
Replace-Text "671×6=" "782×5="
Replace-Text "732×2=" "267×3="

Replace-Text "405×2=" "823×6="
Replace-Text "433×4=" "836×7="
Replace-Text "489×6=" "394×6="
Replace-Text "164×6=" "825×9="
Replace-Text "538×4=" "570×6="
